$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed values for columns B (Ligand symbol) and C (Receptor symbol) for all data rows

$rows = @(
    @("ECs","ECs",3,1,105.2534506666667,315.760352,0.2966477300323703,0.2966477300323703,3,1,37.010628,111.031884,0.9928387522559806,0.9928387522559805,3895.496308340352,35059.46677506317,0.2945233621449075,0.2945233621449074),
    @("ECs","FAPs",3,1,105.2534506666667,315.760352,0.2966477300323703,0.2966477300323703,1,0.3333333333333333,0.05920466666666666,0.177614,0.001588211023629876,0.001588211023629876,6.231495462236444,56.083459160128,0.00047113919497219,0.00047113919497219),
    @("ECs","sCs",3,1,105.2534506666667,315.760352,0.2966477300323703,0.2966477300323703,2,0.6666666666666666,0.2077493333333333,0.623248,0.005573036720389572,0.005573036720389571,21.86633420703289,196.797007863296,0.001653228692490612,0.001653228692490612),
    @("FAPs","ECs",3,1,161.0956266666667,483.28688,0.4540340641196984,0.4540340641196984,3,1,37.010628,111.031884,0.9928387522559806,0.9928387522559805,5962.250310986881,53660.25279888192,0.4507826137023133,0.4507826137023132),
    @("FAPs","FAPs",3,1,161.0956266666667,483.28688,0.4540340641196984,0.4540340641196984,1,0.3333333333333333,0.05920466666666666,0.177614,0.001588211023629876,0.001588211023629876,9.537612878257777,85.83851590431999,0.0007211019057383791,0.0007211019057383791),
    @("FAPs","sCs",3,1,161.0956266666667,483.28688,0.4540340641196984,0.4540340641196984,2,0.6666666666666666,0.2077493333333333,0.623248,0.005573036720389572,0.005573036720389571,33.46750904291556,301.20758138624,0.002530348511646792,0.002530348511646792),
    @("sCs","ECs",3,1,88.46048300000001,265.381449,0.2493182058479313,0.2493182058479313,3,1,37.010628,111.031884,0.9928387522559806,0.9928387522559805,3273.978029013325,29465.80226111992,0.2475327764087598,0.2475327764087598),
    @("sCs","FAPs",3,1,88.46048300000001,265.381449,0.2493182058479313,0.2493182058479313,1,0.3333333333333333,0.05920466666666666,0.177614,0.001588211023629876,0.001588211023629876,5.237273409187334,47.135460682686,0.0003959699229193072,0.0003959699229193072),
    @("sCs","sCs",3,1,88.46048300000001,265.381449,0.2493182058479313,0.2493182058479313,2,0.6666666666666666,0.2077493333333333,0.623248,0.005573036720389572,0.005573036720389571,18.37760636959467,165.398457326352,0.001389459516252167,0.001389459516252167),
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]   # A: Sending cluster
    $ws.Cells.Item($r, 2).Value = "Anxa2"   # B: Ligand symbol
    $ws.Cells.Item($r, 3).Value = "Robo4"   # C: Receptor symbol
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Target cluster
    for ($i = 2; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, 3 + $i).Value = $row[$i]
    }
    $r++
}
